$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = " 18 Sep, 2022"

$ws.Range("B4").Value = "白宇轩"
$ws.Range("E4").Value = "范青桐"
$ws.Range("F4").Value = "邱晨朔"
$ws.Range("H4").Value = "廖从云"
$ws.Range("I4").Value = "边麓元"
$ws.Range("K4").Value = "龚搏扬"

$ws.Range("B5").Value = "卢逸"
$ws.Range("C5").Value = "陈元畅"
$ws.Range("E5").Value = "张宸瑞"
$ws.Range("F5").Value = "曾韦翔"
$ws.Range("H5").Value = "詹悦"
$ws.Range("I5").Value = "林彦含"
$ws.Range("K5").Value = "龙飞宇"

$ws.Range("B6").Value = "骆子墨"
$ws.Range("C6").Value = "赖思轩"
$ws.Range("E6").Value = "程启航"
$ws.Range("F6").Value = "李星宸"
$ws.Range("H6").Value = "黄婧涵"
$ws.Range("I6").Value = "石清泓"
$ws.Range("K6").Value = "迟涵予"
$ws.Range("L6").Value = "郑俊永"

$ws.Range("B7").Value = "杜心扬"
$ws.Range("C7").Value = "章淏博"
$ws.Range("E7").Value = "吴周毅"
$ws.Range("F7").Value = "杨熙宇"
$ws.Range("H7").Value = "张扬"
$ws.Range("I7").Value = "陈李石农"
$ws.Range("K7").Value = "王昊天"
$ws.Range("L7").Value = "丁鹏元"
